$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the existing "sum" header (G1) to the new "Save" header (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new "Save" column values (0) for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
